# NanoServiceAPI Demo.pptx — "Azure Function" slide (sldId 266 / cId 499019303)
#
# The paragraph describing the Hyponatremia function is rewritten:
#   - "The " + "HyponatremiaDDx" (err=1) + " function 1) decides whether
#     it is hyponatremia or not, ..."
#   becomes a single run:
#   "The Hyponatremia Diagnosis Function 1) decides whether there is
#     hyponatremia or not, ..."
# and the text box is widened (cx 11215141 -> 11791283 EMU) while its
# height/position stay the same.

$p = $ppt.ActivePresentation

# sldId 266 is the 12th slide in the deck's sldIdLst.
$s = $p.Slides.Item(12)

# Shape id=9 ("Rectangle 8", creationId {6963E91B-41F0-43F7-A765-18C8AA1A67B0})
# is the 4th shape on the slide.
$sh = $s.Shapes.Item(4)

$tf = $sh.TextFrame
$tr = $tf.TextRange
$para1 = $tr.Paragraphs(1)

$newText = "The Hyponatremia Diagnosis Function 1) decides whether there is hyponatremia or not, 2) adds a new variable “hyponatremia” to the patient actor in case the variable is not yet added, and 3) sets the value to true (hyponatremia) or false (not hyponatremia)."

# Put the full replacement text into the first physical run ...
$r1 = $para1.Runs(1)
$r1.Text = $newText

# ... then clear out the other two physical runs (from the end first, so
# the earlier indices don't shift under us) so the paragraph collapses
# back down to a single <a:r>.
$r3 = $para1.Runs(3)
$r3.Text = ""
$r2 = $para1.Runs(2)
$r2.Text = ""

# Widen the text box; height/position are unchanged.
$sh.Width = 928.4474803149607
